$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (task #3): Status changed from "in progress" to "done"
$ws.Range("D4").Value = "done"

# Row 10 (task #9, "Write a SRS document."): Responsible Unit + Comment filled in
$ws.Range("C10").Value = "all"
$ws.Range("E10").Value = "hbahnev, PeUzunov"

# Row 11 (task #10): fix task wording "Testplan" -> "Test plan"
$ws.Range("B11").Value = "Write a Test plan for functionality test."
